$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '92.053.86'
$ws.Range("E2").Value = '  +1.85%  '

# Row 3: 'Ethereum'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.179.54'
$ws.Range("E3").Value = '  +2.84%  '

# Row 4: 'TetherUSD'
$ws.Range("E4").Value = '  +0.13%  '

# Row 5: 'Solana'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.55'
$ws.Range("E5").Value = '  +3.08%  '

# Row 6: 'BNB'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '623.34'
$ws.Range("E6").Value = '  +0.02%  '

# Row 7: 'XRP'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.13'
$ws.Range("E7").Value = '  +4.91%  '

# Row 8: 'Dogecoin'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.373'
$ws.Range("E8").Value = '  +2.29%  '

# Row 10: 'LidoStakedEther'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.177.22'
$ws.Range("E10").Value = '  +2.82%  '

# Row 11: 'Cardano'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.750'
$ws.Range("E11").Value = '  +3.11%  '

# Row 12: 'TRON'
$ws.Range("E12").Value = '  +3.98%  '

# Row 13: 'ShibaInu'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("E13").Value = '  -0.15%  '

# Row 14: 'Avalanche'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.98'
$ws.Range("E14").Value = '  -0.33%  '

# Row 15: 'Toncoin'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.54'
$ws.Range("E15").Value = '  +1.17%  '

# Row 16: 'WrappedBTC'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.570.05'
$ws.Range("E16").Value = '  +1.72%  '

# Row 17: 'WrappedliquidstakedEther2.0'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.734.31'
$ws.Range("E17").Value = '  +1.98%  '

# Row 18: 'WrappedEther'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.142.07'
$ws.Range("E18").Value = '  +2.17%  '

# Row 19: 'SuiNetwork'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.75'
$ws.Range("E19").Value = '  -3.07%  '

# Row 20: 'Chainlink'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.70'
$ws.Range("E20").Value = '  +12.00%  '

# Row 21: 'PEPE'
$ws.Range("E21").Value = '  -2.14%  '

# Row 22: 'Polkadot'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.87'
$ws.Range("E22").Value = '  +5.15%  '

# Row 23: 'BitcoinCash'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '447.61'
$ws.Range("E23").Value = '  +2.64%  '

# Row 24: 'Uniswap'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.32'
$ws.Range("E24").Value = '  +4.74%  '

# Row 25: 'NEARProtocol'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.27'
$ws.Range("E25").Value = '  +5.60%  '

# Row 26: 'Litecoin'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '90.06'
$ws.Range("E26").Value = '  +1.30%  '

# Row 27: 'Aptos'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.14'
$ws.Range("E27").Value = '  +0.34%  '

# Row 28: 'WrappedeETH'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.271.34'
$ws.Range("E28").Value = '  +0.55%  '

# Row 29: 'Dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.05%  '

# Row 30: 'Hedera'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.140'
$ws.Range("E30").Value = '  +57.23%  '

# Row 31: 'Stellar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.236'
$ws.Range("E31").Value = '  +19.18%  '

# Row 32: 'Cronos'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.174'
$ws.Range("E32").Value = '  +9.65%  '

# Row 33: 'InternetComputer(DFINITY)'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.40'
$ws.Range("E33").Value = '  +0.60%  '

# Row 34: 'Kaspa'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.166'
$ws.Range("E34").Value = '  +6.70%  '

# Row 35: 'Binance-PegBSC-USD' -> 'RenderToken'
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.00'
$ws.Range("E35").Value = '  +11.20%  '

# Row 36: 'RenderToken' -> 'EthereumClassic'
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.79'
$ws.Range("E36").Value = '  +3.81%  '

# Row 37: 'EthereumClassic' -> 'MantraDAO'
$ws.Range("B37").Value = 'MantraDAO'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.19'
$ws.Range("E37").Value = '  +22.32%  '

# Row 38: 'MantraDAO' -> 'Bittensor'
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '521.19'
$ws.Range("E38").Value = '  +3.43%  '

# Row 39: 'Bittensor' -> 'Fetch.AI'
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.36'
$ws.Range("E39").Value = '  +5.32%  '

# Row 40: 'PancakeSwap'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.93'
$ws.Range("E40").Value = '  +1.38%  '

# Row 41: 'Fetch.AI' -> 'PolygonEcosystemToken'
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.454'
$ws.Range("E41").Value = '  +12.08%  '

# Row 43: 'PolygonEcosystemToken' -> 'Binance-PegBSC-USD'
$ws.Range("B43").Value = 'Binance-PegBSC-USD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.809'
$ws.Range("E43").Value = '  -18.83%  '

# Row 44: 'WhiteBITCoin'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.19'
$ws.Range("E44").Value = '  +0.15%  '

# Row 45: 'USDe'
$ws.Range("E45").Value = '  -0.01%  '

# Row 46: 'ARBITRUM'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.722'
$ws.Range("E46").Value = '  +4.17%  '

# Row 47: 'Stacks'
$ws.Range("E47").Value = '  +2.33%  '

# Row 48: 'Monero'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '155.39'
$ws.Range("E48").Value = '  +1.74%  '

# Row 49: 'ImmutableX'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.40'
$ws.Range("E49").Value = '  +3.93%  '

# Row 50: 'Filecoin'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.56'
$ws.Range("E50").Value = '  +2.86%  '

# Row 51: 'OKB' -> 'VeChain'
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0330'
$ws.Range("E51").Value = '  +13.53%  '
